$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update EMILIA's Saldo (row 2, column C): 126143 -> 149741.06
$ws.Cells.Item(2, 3).Value = 149741.06

# 2) Remove the ADRIANA row (005715733 / ADRIANA / 2000.09).
#    In the current (pre-insert) layout this is row 8; removing it now -
#    before the inserts below - keeps the row-5 insertion point unaffected
#    (row 8 is below row 5).
$ws.Rows.Item(8).Delete()

# 3) Insert two new rows immediately before AHMAD's row (row 5) and fill
#    them with the FERNANDA and PRISCILLA accounts. Account numbers keep
#    their leading zeros, so format column A as text before writing them.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "000806386"
$ws.Cells.Item(5, 2).Value = "FERNANDA"
$ws.Cells.Item(5, 3).Value = 49952.51

$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "004224284"
$ws.Cells.Item(6, 2).Value = "PRISCILLA"
$ws.Cells.Item(6, 3).Value = 28903.8

# 4) The old JULIANA row (004813088 / JULIANA / 24295.92) has shifted from
#    row 6 to row 8 after the two inserts above; replace its contents with
#    the VERANICE account.
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "005009947"
$ws.Cells.Item(8, 2).Value = "VERANICE"
$ws.Cells.Item(8, 3).Value = 14952.35
